$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.328.92"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.873.86"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'235.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "'0.2872"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.06591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'21.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07925"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'96.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "1.869.20"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "'0.6933"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "'5.114"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'269.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "30.304.83"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "'14.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "'0.000007722"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "2.118.05"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'5.271"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "'6.213"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'9.412"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").Value = "'167.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'18.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "'1.358"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").Value = "'0.09899"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'4.359"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'1.464"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "'0.04755"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'0.7044"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "'2.724"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").Value = "'0.01874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").Value = "'2.803"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.43%  "
$ws.Range("D40").Value = "'6.199"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").Value = "'72.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "'1.959"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "'0.4184"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "'0.8423"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'102.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "'7.136"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'937.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").Value = "'34.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'0.05692"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.69%  "
